$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
